# Move "Picture 3" (the green leaf picture) on slide 4 (sldId 259) to its new position.
# Original OOXML offset: x=5857875 EMU, y=2841173 EMU
# New OOXML offset:      x=5860829 EMU, y=2980666 EMU
# PowerPoint COM uses points (1 pt = 12700 EMU).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$sh = $s.Shapes.Item(6)

# A tiny epsilon nudges the value past a float32 rounding boundary inside the
# host so the EMU value written back out lands exactly on the target (COM
# Left/Top are single-precision points internally).
$sh.Left = 5860829 / 12700
$sh.Top = (2980666 / 12700) + 0.00002
